# Update market-price / profit figures across the Leve profit sheets.
# Values reflect a refreshed price snapshot from the scheduled market-data run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 950
$ws.Range("I49").Value = 500
$ws.Range("K49").Value = 1500
$ws.Range("M49").Value = -1364
$ws.Range("H62").Value = 72391.914
$ws.Range("I62").Value = 102587.875
$ws.Range("K62").Value = 102587.875
$ws.Range("M62").Value = -101963.875
$ws.Range("H64").Value = 8437.615
$ws.Range("I64").Value = 3995
$ws.Range("J64").Value = 9245.362999999999
$ws.Range("K64").Value = 3995
$ws.Range("L64").Value = 9245.362999999999
$ws.Range("M64").Value = -3747
$ws.Range("N64").Value = -9741.362999999999
$ws.Range("H65").Value = 72391.914
$ws.Range("I65").Value = 102587.875
$ws.Range("K65").Value = 512939.375
$ws.Range("M65").Value = -509819.375
$ws.Range("H67").Value = 8437.615
$ws.Range("I67").Value = 3995
$ws.Range("J67").Value = 9245.362999999999
$ws.Range("K67").Value = 3995
$ws.Range("L67").Value = 9245.362999999999
$ws.Range("M67").Value = -3137
$ws.Range("N67").Value = -10961.363
$ws.Range("H80").Value = 879.2857
$ws.Range("I80").Value = 494.625
$ws.Range("J80").Value = 1392.1666
$ws.Range("K80").Value = 1483.875
$ws.Range("L80").Value = 4176.4998
$ws.Range("M80").Value = -485.875
$ws.Range("N80").Value = -6172.4998
$ws.Range("H83").Value = 879.2857
$ws.Range("I83").Value = 494.625
$ws.Range("J83").Value = 1392.1666
$ws.Range("K83").Value = 4451.625
$ws.Range("L83").Value = 12529.4994
$ws.Range("M83").Value = 540.375
$ws.Range("N83").Value = -22513.4994
$ws.Range("H100").Value = 6322.4287
$ws.Range("I100").Value = 2224.75
$ws.Range("K100").Value = 2224.75
$ws.Range("M100").Value = -1683.75
$ws.Range("H131").Value = 9898.666999999999
$ws.Range("I131").Value = 9023.5
$ws.Range("K131").Value = 27070.5
$ws.Range("M131").Value = -22030.5
$ws.Range("H139").Value = 70416.5
$ws.Range("J139").Value = 70416.5
$ws.Range("L139").Value = 70416.5
$ws.Range("N139").Value = -80696.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 491.5
$ws.Range("I2").Value = 527.25
$ws.Range("J2").Value = 420
$ws.Range("K2").Value = 527.25
$ws.Range("L2").Value = 420
$ws.Range("M2").Value = -414.25
$ws.Range("N2").Value = -646
$ws.Range("H22").Value = 4411.067
$ws.Range("J22").Value = 15000
$ws.Range("L22").Value = 15000
$ws.Range("N22").Value = -15598
$ws.Range("H37").Value = 26250.25
$ws.Range("J37").Value = 40001
$ws.Range("L37").Value = 40001
$ws.Range("N37").Value = -40547
$ws.Range("H97").Value = 2439.6667
$ws.Range("I97").Value = 2439.6667
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2439.6667
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1943.6667
$ws.Range("N97").Value = ""
$ws.Range("H102").Value = 2995.353
$ws.Range("I102").Value = 2061.4
$ws.Range("K102").Value = 2061.4
$ws.Range("M102").Value = -439.4000000000001
$ws.Range("H116").Value = 491.5
$ws.Range("I116").Value = 527.25
$ws.Range("J116").Value = 420
$ws.Range("K116").Value = 527.25
$ws.Range("L116").Value = 420
$ws.Range("M116").Value = 1766.75
$ws.Range("N116").Value = -5008

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 491.5
$ws.Range("I3").Value = 527.25
$ws.Range("J3").Value = 420
$ws.Range("K3").Value = 527.25
$ws.Range("L3").Value = 420
$ws.Range("M3").Value = -413.25
$ws.Range("N3").Value = -648
$ws.Range("H22").Value = 812.2
$ws.Range("I22").Value = 598.875
$ws.Range("J22").Value = 1665.5
$ws.Range("K22").Value = 598.875
$ws.Range("L22").Value = 1665.5
$ws.Range("M22").Value = -425.875
$ws.Range("N22").Value = -2011.5
$ws.Range("H96").Value = 22030.125
$ws.Range("I96").Value = 6874.3335
$ws.Range("K96").Value = 6874.3335
$ws.Range("M96").Value = -4128.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4471.4
$ws.Range("I16").Value = 4363.875
$ws.Range("J16").Value = 4901.5
$ws.Range("K16").Value = 4363.875
$ws.Range("L16").Value = 4901.5
$ws.Range("M16").Value = -4076.875
$ws.Range("N16").Value = -5475.5
$ws.Range("H59").Value = 34713.285
$ws.Range("J59").Value = 39818.6
$ws.Range("L59").Value = 39818.6
$ws.Range("N59").Value = -42108.6
$ws.Range("H113").Value = 4471.4
$ws.Range("I113").Value = 4363.875
$ws.Range("J113").Value = 4901.5
$ws.Range("K113").Value = 4363.875
$ws.Range("L113").Value = 4901.5
$ws.Range("M113").Value = -2193.875
$ws.Range("N113").Value = -9241.5
$ws.Range("H134").Value = 2940.3794
$ws.Range("I134").Value = 2798.0476
$ws.Range("J134").Value = 3314
$ws.Range("K134").Value = 8394.1428
$ws.Range("L134").Value = 9942
$ws.Range("M134").Value = -5859.1428
$ws.Range("N134").Value = -15012

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1000
$ws.Range("I63").Value = 1000
$ws.Range("K63").Value = 3000
$ws.Range("M63").Value = -2251
$ws.Range("H64").Value = 1349.5
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""
$ws.Range("H66").Value = 1000
$ws.Range("I66").Value = 1000
$ws.Range("K66").Value = 9000
$ws.Range("M66").Value = -5256
$ws.Range("H67").Value = 1349.5
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""
$ws.Range("H76").Value = 1000
$ws.Range("I76").Value = 1000
$ws.Range("K76").Value = 3000
$ws.Range("M76").Value = -2617
$ws.Range("H79").Value = 1000
$ws.Range("I79").Value = 1000
$ws.Range("K79").Value = 3000
$ws.Range("M79").Value = -1674
$ws.Range("H87").Value = 4014
$ws.Range("I87").Value = 4014
$ws.Range("K87").Value = 12042
$ws.Range("M87").Value = -10794
$ws.Range("H90").Value = 4014
$ws.Range("I90").Value = 4014
$ws.Range("K90").Value = 36126
$ws.Range("M90").Value = -29886
$ws.Range("H117").Value = 1803.0555
$ws.Range("I117").Value = 2509.8333
$ws.Range("J117").Value = 1449.6666
$ws.Range("K117").Value = 7529.499899999999
$ws.Range("L117").Value = 4348.9998
$ws.Range("M117").Value = -4087.499899999999
$ws.Range("N117").Value = -11232.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = ""
$ws.Range("H122").Value = 12273.704
$ws.Range("I122").Value = 13642.19
$ws.Range("J122").Value = 7484
$ws.Range("K122").Value = 40926.57
$ws.Range("L122").Value = 22452
$ws.Range("M122").Value = -38476.57
$ws.Range("N122").Value = -27352
$ws.Range("H126").Value = 3007.7144
$ws.Range("I126").Value = 2993.1667
$ws.Range("K126").Value = 8979.500100000001
$ws.Range("M126").Value = -6509.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7037.9
$ws.Range("I61").Value = 4800
$ws.Range("J61").Value = 7997
$ws.Range("K61").Value = 4800
$ws.Range("L61").Value = 7997
$ws.Range("M61").Value = -4598
$ws.Range("N61").Value = -8401
$ws.Range("H113").Value = 7037.9
$ws.Range("I113").Value = 4800
$ws.Range("J113").Value = 7997
$ws.Range("K113").Value = 4800
$ws.Range("L113").Value = 7997
$ws.Range("M113").Value = -2630
$ws.Range("N113").Value = -12337
$ws.Range("H127").Value = 185399.5
$ws.Range("J127").Value = 185399.5
$ws.Range("L127").Value = 185399.5
$ws.Range("N127").Value = -195319.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2593943
$ws.Range("I2").Value = 2593943
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2593943
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2593831
$ws.Range("N2").Value = ""
$ws.Range("H3").Value = 200
$ws.Range("J3").Value = 200
$ws.Range("L3").Value = 200
$ws.Range("N3").Value = -428
$ws.Range("H8").Value = 9999
$ws.Range("I8").Value = 9998
$ws.Range("K8").Value = 9998
$ws.Range("M8").Value = -9858
$ws.Range("H62").Value = 14324.2
$ws.Range("I62").Value = 4498.6665
$ws.Range("K62").Value = 4498.6665
$ws.Range("M62").Value = -3874.6665
$ws.Range("H65").Value = 14324.2
$ws.Range("I65").Value = 4498.6665
$ws.Range("K65").Value = 22493.3325
$ws.Range("M65").Value = -19373.3325
$ws.Range("H107").Value = 1958.45
$ws.Range("I107").Value = 2291.4
$ws.Range("J107").Value = 959.6
$ws.Range("K107").Value = 6874.200000000001
$ws.Range("L107").Value = 2878.8
$ws.Range("M107").Value = -4954.200000000001
$ws.Range("N107").Value = -6718.8
$ws.Range("H122").Value = 5195.16
$ws.Range("I122").Value = 3305
$ws.Range("K122").Value = 9915
$ws.Range("M122").Value = -7465
$ws.Range("H126").Value = 1768.25
$ws.Range("I126").Value = 1646.1428
$ws.Range("K126").Value = 4938.428400000001
$ws.Range("M126").Value = -2468.428400000001
$ws.Range("H132").Value = 1964.3715
$ws.Range("I132").Value = 1560.4138
$ws.Range("J132").Value = 3916.8333
$ws.Range("K132").Value = 4681.2414
$ws.Range("L132").Value = 11750.4999
$ws.Range("M132").Value = -2151.2414
$ws.Range("N132").Value = -16810.4999
